$d = $word.ActiveDocument

# Locate the paragraph containing the "Second, flooding is tremendously
# robust. ..." sentence (it is currently split across three runs:
# "...are ", "literally ", "blown...their ").
$hit = $d.Content
$hit.Find.Execute("Second, flooding is tremendously robust*know their neighbors.",
    $false, $false, $true, $false, $false, $true, 1, $false, "", 0)

$p = $hit.Paragraphs(1)
$r = $p.Range

# Merge the "are " / "literally " / "blown..." runs into a single run by
# replacing the text that spans the run boundary with itself (on a
# duplicate range, so $r keeps covering the whole paragraph) - Word
# coalesces runs that end up with identical formatting.
$mergeRange = $r.Duplicate
$mergeRange.Find.Execute("are literally blown", $true, $false, $false, $false, $false,
    $true, 1, $false, "are literally blown", 2)

# Bump the paragraph mark and every run in the paragraph to 12pt (sz/szCs
# 24 half-points), matching the surrounding paragraphs' explicit sizing.
$r.Font.Size = 12
$r.Font.SizeBi = 12
